$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style A2 fully (bold font, thin box border, centered/top aligned)
$r = $ws.Range("A2")
$r.Font.Bold = $true
$r.Borders.LineStyle = 1
$r.Borders.Weight = 2
$r.HorizontalAlignment = -4108
$r.VerticalAlignment = -4160

# Clone the exact same style onto B1 in a single atomic operation so no
# intermediate/orphan cell-format entries are produced.
$r.Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
